$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (C) column date for rows 2-7 from 2023-10-05 (45204)
# to 2023-10-08 (45207), keeping the existing date formatting.
$newDate = Get-Date -Year 2023 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0

foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
